$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.007.90'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.58%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.171.72'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.37%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.72'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.610'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.65'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -5.47%  '

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.573'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -6.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.72'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0925'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.101'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.21%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.75'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.498.96'
$ws.Range("D14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.86'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.809'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.160.85'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '40.940.38'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000101'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -8.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.35'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.94'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.87'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '224.98'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.95'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -8.03%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.91'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -6.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.55'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.76%  '

$ws.Range("E29").Value = '  -1.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.14'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.80'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.13%  '

$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '30.63'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.09%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0774'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.13'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -8.17%  '

$ws.Range("E35").Value = '  -3.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.103'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -9.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.14'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0288'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.40'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.05'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.43'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '59.79'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.58%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.189'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.31'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0973'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '98.89'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.08'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.53%  '

$ws.Range("E48").Value = '  -3.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.22'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -7.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.65'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.29%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.377.63'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.33%  '
